$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 6151
$ws.Cells.Item(3, 7).Value = 8320
$ws.Cells.Item(3, 10).Value = 6543
$ws.Cells.Item(4, 10).Value = 1419
$ws.Cells.Item(5, 10).Value = 501
$ws.Cells.Item(6, 10).Value = 8508
$ws.Cells.Item(7, 7).Value = 24699
$ws.Cells.Item(7, 10).Value = 23122

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 185
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(6, 10).Value = 174
$ws.Cells.Item(7, 10).Value = 681
$ws.Cells.Item(8, 10).Value = 1457
$ws.Cells.Item(9, 10).Value = 117
$ws.Cells.Item(11, 10).Value = 371
$ws.Cells.Item(12, 10).Value = 45
$ws.Cells.Item(15, 10).Value = 258
$ws.Cells.Item(16, 10).Value = 91
$ws.Cells.Item(17, 10).Value = 32
$ws.Cells.Item(19, 10).Value = 686
$ws.Cells.Item(27, 10).Value = 142
$ws.Cells.Item(29, 10).Value = 1275
$ws.Cells.Item(33, 10).Value = 1058
$ws.Cells.Item(34, 10).Value = 105
$ws.Cells.Item(36, 10).Value = 314
$ws.Cells.Item(42, 7).Value = 995
$ws.Cells.Item(42, 10).Value = 977
$ws.Cells.Item(43, 10).Value = 197
$ws.Cells.Item(46, 10).Value = 76
$ws.Cells.Item(49, 10).Value = 153
$ws.Cells.Item(50, 10).Value = 140
$ws.Cells.Item(51, 10).Value = 291
$ws.Cells.Item(52, 10).Value = 579
$ws.Cells.Item(53, 10).Value = 324
$ws.Cells.Item(54, 10).Value = 446
$ws.Cells.Item(55, 10).Value = 331
$ws.Cells.Item(57, 10).Value = 103
$ws.Cells.Item(62, 10).Value = 8
$ws.Cells.Item(63, 10).Value = 84
$ws.Cells.Item(64, 10).Value = 152
$ws.Cells.Item(65, 10).Value = 571
$ws.Cells.Item(67, 10).Value = 875
$ws.Cells.Item(68, 10).Value = 48
$ws.Cells.Item(70, 10).Value = 32
$ws.Cells.Item(73, 10).Value = 222
$ws.Cells.Item(76, 10).Value = 353
$ws.Cells.Item(78, 10).Value = 278
$ws.Cells.Item(79, 10).Value = 657
$ws.Cells.Item(80, 10).Value = 37
$ws.Cells.Item(83, 10).Value = 462
$ws.Cells.Item(84, 10).Value = 194
$ws.Cells.Item(85, 10).Value = 950
$ws.Cells.Item(90, 10).Value = 250
$ws.Cells.Item(93, 10).Value = 101
$ws.Cells.Item(94, 10).Value = 240
$ws.Cells.Item(98, 10).Value = 168
$ws.Cells.Item(100, 10).Value = 43
$ws.Cells.Item(101, 7).Value = 24699
$ws.Cells.Item(101, 10).Value = 23122

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 211
$ws.Cells.Item(6, 10).Value = 219
$ws.Cells.Item(7, 10).Value = 681

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 10).Value = 110
$ws.Cells.Item(7, 10).Value = 371

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 252
$ws.Cells.Item(3, 10).Value = 336
$ws.Cells.Item(4, 10).Value = 65
$ws.Cells.Item(6, 10).Value = 277
$ws.Cells.Item(7, 10).Value = 950

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 10).Value = 137
$ws.Cells.Item(6, 10).Value = 236
$ws.Cells.Item(7, 10).Value = 579

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 10).Value = 214
$ws.Cells.Item(7, 10).Value = 324

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 397
$ws.Cells.Item(3, 10).Value = 438
$ws.Cells.Item(6, 10).Value = 508
$ws.Cells.Item(7, 10).Value = 1457

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 174
$ws.Cells.Item(4, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 462

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 249
$ws.Cells.Item(3, 10).Value = 352
$ws.Cells.Item(5, 10).Value = 44
$ws.Cells.Item(7, 10).Value = 1058

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 10).Value = 168
$ws.Cells.Item(6, 10).Value = 201
$ws.Cells.Item(7, 10).Value = 571

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 10).Value = 218
$ws.Cells.Item(7, 10).Value = 875

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 10).Value = 64
$ws.Cells.Item(6, 10).Value = 59
$ws.Cells.Item(7, 10).Value = 194

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(3, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 153

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 10).Value = 213
$ws.Cells.Item(7, 10).Value = 446

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 392
$ws.Cells.Item(3, 10).Value = 441
$ws.Cells.Item(6, 10).Value = 325
$ws.Cells.Item(7, 10).Value = 1275

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(6, 10).Value = 265
$ws.Cells.Item(7, 10).Value = 686

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 10).Value = 196
$ws.Cells.Item(7, 10).Value = 353

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(2, 10).Value = 52
$ws.Cells.Item(6, 10).Value = 65
$ws.Cells.Item(7, 10).Value = 174

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 211
$ws.Cells.Item(3, 7).Value = 360
$ws.Cells.Item(7, 7).Value = 995
$ws.Cells.Item(7, 10).Value = 977

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(7, 10).Value = 278

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 10).Value = 174
$ws.Cells.Item(7, 10).Value = 331

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 10).Value = 17
$ws.Cells.Item(7, 10).Value = 76

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(6, 10).Value = 194
$ws.Cells.Item(7, 10).Value = 657

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 10).Value = 42
$ws.Cells.Item(7, 10).Value = 152

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 314

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(3, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 101

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(6, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 43

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 10).Value = 39
$ws.Cells.Item(7, 10).Value = 105

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 10).Value = 134
$ws.Cells.Item(7, 10).Value = 240

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(4, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 258

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(6, 10).Value = 104
$ws.Cells.Item(7, 10).Value = 168

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 10).Value = 46
$ws.Cells.Item(7, 10).Value = 140

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(3, 10).Value = 39
$ws.Cells.Item(7, 10).Value = 117

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 10).Value = 74
$ws.Cells.Item(7, 10).Value = 222

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 10).Value = 71
$ws.Cells.Item(7, 10).Value = 185

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(3, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(6, 10).Value = 49
$ws.Cells.Item(7, 10).Value = 142

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 10).Value = 70
$ws.Cells.Item(7, 10).Value = 250

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(4, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 291

$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(6, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 48

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 103

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(7, 10).Value = 197

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 37

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(2, 10).Value = 34
$ws.Cells.Item(7, 10).Value = 100

$ws = $wb.Worksheets.Item('Beverly')
$ws.Cells.Item(3, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 45

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(2, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 91

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(7, 10).Value = 8
